$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings) - force text so Excel does not
# auto-convert the "<Month> <Year>" strings into date serials.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "February 2025"
$ws.Range("A1").Style = "Normal"

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "March 2025"
$ws.Range("G1").Style = "Normal"

# Update data row values
$ws.Range("A2").Value = 1.459
$ws.Range("B2").Value = 0.225
$ws.Range("C2").Value = 0.118
$ws.Range("D2").Value = 0.295
$ws.Range("E2").Value = 0.009
$ws.Range("F2").Value = -0.329
$ws.Range("G2").Value = 1.778
